$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins/Losses/Ties, matching the style of the
# existing header row (e.g. AB1) so they get the bold/centered/bordered look.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 29).Value = 84
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 1
}
